# Informe Tecnico.docx - "commit 31 mayo posible final"
# - Corrects the memo date from 29 May to 31 May 2018.
# - Normalizes the signature-block names/titles (proper case, and adds
#   the missing space after "Cargo:").

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replacement not found -> $find"
    }
}

# 1. FECHA:  29 May de 2018  ->  31 May de 2018
Replace-Text "29 May de 2018" "31 May de 2018"

# 2. Signature table, left column - requester name
Replace-Text "Nombre: ECON. ALEXIA DUQUE" "Nombre: Econ. Alexia Duque"

# 3. Signature table, right column - approver name
Replace-Text "Nombre: ING. ANA MARIA ROMERO MENDOZA" "Nombre: Ing. Ana Maria Romero Mendoza"

# 4. Signature table, left column - requester title (also adds the
#    missing space after the colon)
Replace-Text "Cargo:COORDINADORA DE SERVICIOS GENERALES" "Cargo: Coordinadora De Servicios Generales"

# 5. Signature table, right column - approver title (also adds the
#    missing space after the colon)
Replace-Text "Cargo:Directora Administrativa" "Cargo: Directora Administrativa"
